$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 72.79331566666666
$ws.Range("H2").Value = 218.379947
$ws.Range("I2").Value = 0.2828741606141505
$ws.Range("J2").Value = 0.2828741606141506
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.16653
$ws.Range("N2").Value = 51.49959
$ws.Range("O2").Value = 0.0560345397128279
$ws.Range("P2").Value = 0.0560345397128279
$ws.Range("Q2").Value = 1249.608637191303
$ws.Range("R2").Value = 11246.47773472173
$ws.Range("S2").Value = 0.01585072338666647
$ws.Range("T2").Value = 0.01585072338666648

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 72.79331566666666
$ws.Range("H3").Value = 218.379947
$ws.Range("I3").Value = 0.2828741606141505
$ws.Range("J3").Value = 0.2828741606141506
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 256.4443053333333
$ws.Range("N3").Value = 769.332916
$ws.Range("O3").Value = 0.8370788162388805
$ws.Range("P3").Value = 0.8370788162388805
$ws.Range("Q3").Value = 18667.43126904838
$ws.Range("R3").Value = 168006.8814214354
$ws.Range("S3").Value = 0.2367879675114601
$ws.Range("T3").Value = 0.2367879675114601

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 72.79331566666666
$ws.Range("H4").Value = 218.379947
$ws.Range("I4").Value = 0.2828741606141505
$ws.Range("J4").Value = 0.2828741606141506
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.74538866666666
$ws.Range("N4").Value = 98.236166
$ws.Range("O4").Value = 0.1068866440482915
$ws.Range("P4").Value = 0.1068866440482915
$ws.Range("Q4").Value = 2383.645413840355
$ws.Range("R4").Value = 21452.8087245632
$ws.Range("S4").Value = 0.03023546971602396
$ws.Range("T4").Value = 0.03023546971602396

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 117.1700846666667
$ws.Range("H5").Value = 351.510254
$ws.Range("I5").Value = 0.4553218801152877
$ws.Range("J5").Value = 0.4553218801152878
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.16653
$ws.Range("N5").Value = 51.49959
$ws.Range("O5").Value = 0.0560345397128279
$ws.Range("P5").Value = 0.0560345397128279
$ws.Range("Q5").Value = 2011.403773532873
$ws.Range("R5").Value = 18102.63396179586
$ws.Range("S5").Value = 0.02551375197343955
$ws.Range("T5").Value = 0.02551375197343956

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 117.1700846666667
$ws.Range("H6").Value = 351.510254
$ws.Range("I6").Value = 0.4553218801152877
$ws.Range("J6").Value = 0.4553218801152878
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 256.4443053333333
$ws.Range("N6").Value = 769.332916
$ws.Range("O6").Value = 0.8370788162388805
$ws.Range("P6").Value = 0.8370788162388805
$ws.Range("Q6").Value = 30047.60096819119
$ws.Range("R6").Value = 270428.4087137207
$ws.Range("S6").Value = 0.3811403004145665
$ws.Range("T6").Value = 0.3811403004145665

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 117.1700846666667
$ws.Range("H7").Value = 351.510254
$ws.Range("I7").Value = 0.4553218801152877
$ws.Range("J7").Value = 0.4553218801152878
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 32.74538866666666
$ws.Range("N7").Value = 98.236166
$ws.Range("O7").Value = 0.1068866440482915
$ws.Range("P7").Value = 0.1068866440482915
$ws.Range("Q7").Value = 3836.779962516241
$ws.Range("R7").Value = 34531.01966264617
$ws.Range("S7").Value = 0.04866782772728162
$ws.Range("T7").Value = 0.04866782772728163

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 67.37122333333333
$ws.Range("H8").Value = 202.11367
$ws.Range("I8").Value = 0.2618039592705617
$ws.Range("J8").Value = 0.2618039592705618
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 17.16653
$ws.Range("N8").Value = 51.49959
$ws.Range("O8").Value = 0.0560345397128279
$ws.Range("P8").Value = 0.0560345397128279
$ws.Range("Q8").Value = 1156.530126488366
$ws.Range("R8").Value = 10408.7711383953
$ws.Range("S8").Value = 0.01467006435272187
$ws.Range("T8").Value = 0.01467006435272187

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 67.37122333333333
$ws.Range("H9").Value = 202.11367
$ws.Range("I9").Value = 0.2618039592705617
$ws.Range("J9").Value = 0.2618039592705618
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 256.4443053333333
$ws.Range("N9").Value = 769.332916
$ws.Range("O9").Value = 0.8370788162388805
$ws.Range("P9").Value = 0.8370788162388805
$ws.Range("Q9").Value = 17276.96656717352
$ws.Range("R9").Value = 155492.6991045617
$ws.Range("S9").Value = 0.2191505483128539
$ws.Range("T9").Value = 0.2191505483128539

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 67.37122333333333
$ws.Range("H10").Value = 202.11367
$ws.Range("I10").Value = 0.2618039592705617
$ws.Range("J10").Value = 0.2618039592705618
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 32.74538866666666
$ws.Range("N10").Value = 98.236166
$ws.Range("O10").Value = 0.1068866440482915
$ws.Range("P10").Value = 0.1068866440482915
$ws.Range("Q10").Value = 2206.096892998802
$ws.Range("R10").Value = 19854.87203698922
$ws.Range("S10").Value = 0.02798334660498594
$ws.Range("T10").Value = 0.02798334660498595

